# ---------------------------------------------------------------------------
# Applies the "added partial training of fft model" commit:
#   1. Adds an empty "_GoBack" bookmark right at the very start of the
#      document (inside the Heading1 "Introduction" paragraph, before the
#      "Introduction" run).
#   2. Merges the run-fragments of the "Eventually, this automated ..."
#      sentence (which were split on ~line-width boundaries) into one run,
#      leaving the neighbouring "methodologies" run intact.
#   3. Removes the old "_GoBack" bookmark that used to sit in the middle of
#      the "... discussion of conclusions / of the results and ..." sentence,
#      and merges the now-adjacent " " and "of the results and " runs.
#   4. Merges the three runs of the "Code of Federal Regulations ..."
#      reference entry into a single run.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 0) The document already has one (hidden) "_GoBack" bookmark sitting in the
#    "... discussion of conclusions [of the results and] ..." sentence. This
#    host does not enforce bookmark-name uniqueness the way real Word does
#    (adding another bookmark called "_GoBack" creates a *second* one rather
#    than relocating the existing one, and Bookmarks.Item("_GoBack") then
#    resolves to whichever one, ambiguously) so get rid of the original
#    first and remember where it was -- before minting the new one used in
#    step 1 below.
# ---------------------------------------------------------------------------
$oldBm = $d.Bookmarks.Item("_GoBack")
$oldBmStart = $oldBm.Start
$oldBm.Delete()

# ---------------------------------------------------------------------------
# 1) Insert the "_GoBack" bookmark at the very beginning of the document.
#
#    NOTE: $d.Range(0, 0) / a Selection collapsed to absolute position 0 is
#    special-cased by this host (it behaves like "whole document" instead of
#    an empty range at offset 0), so Bookmarks.Add can't be pointed at it
#    directly -- doing so silently anchors bookmarkEnd to the *next*
#    paragraph. Work around it by typing a one-character placeholder at the
#    front (landing the bookmark on a non-zero, well-behaved range), then
#    deleting just that character with Range.Delete() (which -- unlike a
#    Find/Replace -- does not trigger the "re-merge same-format runs in this
#    paragraph" behaviour and so doesn't disturb the bookmark markers).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Introduction", $true, $false, $false, $false, $false,
                         $true, 1, $false, "zIntroduction", 2)

$placeholder = $d.Range(0, 1)
$d.Bookmarks.Add("_GoBack", $placeholder)

$d.Range(0, 1).Delete()

# ---------------------------------------------------------------------------
# 2) Merge the "Eventually, this automated ... for further " run fragments.
#
#    A Find/Replace on this host re-coalesces every run sharing the same
#    (empty) formatting across the whole paragraph, which would also swallow
#    the preceding "methodologies" run. Drop a throw-away bookmark right
#    after "methodologies" as a wall so the merge can't creep past it, then
#    remove the wall again afterwards.
# ---------------------------------------------------------------------------
$wallRange = $d.Range(0, 0)
$wallRange.Find.Execute("methodologies", $false, $false, $false, $false,
                         $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("ZZWALL1", $d.Range($wallRange.End, $wallRange.End))

$sentence = ". Eventually, this automated system would remove the need for the multiple work stations and staff providing significant cost reductions for the community noise fly-over capability. The automated system also would provide increased accuracy and consistency of the classification thus increased efficacy of the test for further "
$d.Content.Find.Execute($sentence, $true, $false, $false, $false, $false,
                         $true, 1, $false, $sentence, 2)

$d.Bookmarks.Item("ZZWALL1").Delete()

# ---------------------------------------------------------------------------
# 3) Merge the space-run with the "of the results and " run (at the spot the
#    old "_GoBack" bookmark -- removed in step 0 above -- used to occupy),
#    without pulling in the neighbouring "conclusions" / "recommendations
#    ..." runs.
# ---------------------------------------------------------------------------
$d.Bookmarks.Add("ZZWALL2", $d.Range($oldBmStart - 1, $oldBmStart - 1))
$d.Bookmarks.Add("ZZWALL3", $d.Range($oldBmStart + 19, $oldBmStart + 19))

$d.Content.Find.Execute(" of the results and ", $true, $false, $false, $false,
                         $false, $true, 1, $false, " of the results and ", 2)

$d.Bookmarks.Item("ZZWALL2").Delete()
$d.Bookmarks.Item("ZZWALL3").Delete()

# ---------------------------------------------------------------------------
# 4) Merge the "Code of Federal Regulations ... Administration" reference
#    runs into a single run. This paragraph is isolated (its own <w:p>) so
#    no wall is required.
# ---------------------------------------------------------------------------
$reference = "Code of Federal Regulations. (2016). Title 14, Part 36, Noise Standards: Aircraft Type and Airworthiness Certification. Washington, D.C.: Federal Aviation Administration."
$d.Content.Find.Execute($reference, $true, $false, $false, $false, $false,
                         $true, 1, $false, $reference, 2)
